# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.442.49"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.725.98"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.20"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4862"
$ws.Range("E7").Value = "  +1.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2611"
$ws.Range("E8").Value = "  -2.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06190"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.732.68"
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07009"
$ws.Range("E11").Value = "  -2.17%  "
$ws.Range("E12").Value = "  -1.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.535"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5991"
$ws.Range("E14").Value = "  -2.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.25"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.453.19"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007177"
$ws.Range("E19").Value = "  +3.19%  "
$ws.Range("E20").Value = "  -2.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.953.26"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.495"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.581"
$ws.Range("E23").Value = "  -3.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.180"
$ws.Range("E24").Value = "  -2.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.39"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("E26").Value = "  -0.73%  "
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "106.95"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.719"
$ws.Range("E29").Value = "  -4.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.958"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07953"
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.692"
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04520"
$ws.Range("E33").Value = "  -1.60%  "
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9992"
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6235"
$ws.Range("E36").Value = "  -2.08%  "
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.996"
$ws.Range("E38").Value = "  -4.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.398"
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01487"
$ws.Range("E41").Value = "  -1.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.06"
$ws.Range("E42").Value = "  -4.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.399"
$ws.Range("E43").Value = "  -3.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3860"
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.687"
$ws.Range("E45").Value = "  -4.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1154"
$ws.Range("E46").Value = "  -2.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05358"
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.20"
$ws.Range("E48").Value = "  -2.50%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.711"
$ws.Range("E49").Value = "  -2.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.251"
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.01"
$ws.Range("E51").Value = "  -0.64%  "
